$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 3145.365
$ws.Cells.Item(137, 9).Value = 1172.7142
$ws.Cells.Item(137, 10).Value = 7090.6665
$ws.Cells.Item(137, 11).Value = 3518.1426
$ws.Cells.Item(137, 12).Value = 21271.9995
$ws.Cells.Item(137, 13).Value = -968.1425999999997
$ws.Cells.Item(137, 14).Value = -26371.9995

$ws.Cells.Item(138, 8).Value = 2196.6858
$ws.Cells.Item(138, 9).Value = 1531.6666
$ws.Cells.Item(138, 10).Value = 3647.6365
$ws.Cells.Item(138, 11).Value = 4594.9998
$ws.Cells.Item(138, 12).Value = 10942.9095
$ws.Cells.Item(138, 13).Value = 545.0002000000004
$ws.Cells.Item(138, 14).Value = -21222.9095

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4392.25
$ws.Cells.Item(32, 9).Value = 3678.2388
$ws.Cells.Item(32, 11).Value = 3678.2388
$ws.Cells.Item(32, 13).Value = -3391.2388

$ws.Cells.Item(61, 8).Value = 1381.6
$ws.Cells.Item(61, 9).Value = 1005.2308
$ws.Cells.Item(61, 10).Value = 2468.889
$ws.Cells.Item(61, 11).Value = 1005.2308
$ws.Cells.Item(61, 12).Value = 2468.889
$ws.Cells.Item(61, 13).Value = -793.2308
$ws.Cells.Item(61, 14).Value = -2892.889

$ws.Cells.Item(132, 8).Value = 12895.263
$ws.Cells.Item(132, 9).Value = 12182.4
$ws.Cells.Item(132, 10).Value = 13687.333
$ws.Cells.Item(132, 11).Value = 36547.2
$ws.Cells.Item(132, 12).Value = 41061.999
$ws.Cells.Item(132, 13).Value = -34017.2
$ws.Cells.Item(132, 14).Value = -46121.999

$ws.Cells.Item(136, 8).Value = 1381.6
$ws.Cells.Item(136, 9).Value = 1005.2308
$ws.Cells.Item(136, 10).Value = 2468.889
$ws.Cells.Item(136, 11).Value = 3015.6924
$ws.Cells.Item(136, 12).Value = 7406.667
$ws.Cells.Item(136, 13).Value = -465.6923999999999
$ws.Cells.Item(136, 14).Value = -12506.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1207.4
$ws.Cells.Item(134, 9).Value = 938.4583
$ws.Cells.Item(134, 10).Value = 1794.1818
$ws.Cells.Item(134, 11).Value = 2815.3749
$ws.Cells.Item(134, 12).Value = 5382.5454
$ws.Cells.Item(134, 13).Value = -280.3748999999998
$ws.Cells.Item(134, 14).Value = -10452.5454

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 3183.1304
$ws.Cells.Item(86, 9).Value = 2879.3333
$ws.Cells.Item(86, 10).Value = 3378.4285
$ws.Cells.Item(86, 11).Value = 2879.3333
$ws.Cells.Item(86, 12).Value = 3378.4285
$ws.Cells.Item(86, 13).Value = -1756.3333
$ws.Cells.Item(86, 14).Value = -5624.4285

$ws.Cells.Item(89, 8).Value = 3183.1304
$ws.Cells.Item(89, 9).Value = 2879.3333
$ws.Cells.Item(89, 10).Value = 3378.4285
$ws.Cells.Item(89, 11).Value = 14396.6665
$ws.Cells.Item(89, 12).Value = 16892.1425
$ws.Cells.Item(89, 13).Value = -8780.666499999999
$ws.Cells.Item(89, 14).Value = -28124.1425

$ws.Cells.Item(94, 8).Value = 575.6667
$ws.Cells.Item(94, 9).Value = 800
$ws.Cells.Item(94, 10).Value = 530.8
$ws.Cells.Item(94, 11).Value = 800
$ws.Cells.Item(94, 12).Value = 530.8
$ws.Cells.Item(94, 13).Value = -349
$ws.Cells.Item(94, 14).Value = -1432.8

$ws.Cells.Item(99, 8).Value = 3174.1538
$ws.Cells.Item(99, 9).Value = 3227.2
$ws.Cells.Item(99, 10).Value = 3141
$ws.Cells.Item(99, 11).Value = 3227.2
$ws.Cells.Item(99, 12).Value = 3141
$ws.Cells.Item(99, 13).Value = -1729.2
$ws.Cells.Item(99, 14).Value = -6137

$ws.Cells.Item(126, 8).Value = 3174.1538
$ws.Cells.Item(126, 9).Value = 3227.2
$ws.Cells.Item(126, 10).Value = 3141
$ws.Cells.Item(126, 11).Value = 9681.599999999999
$ws.Cells.Item(126, 12).Value = 9423
$ws.Cells.Item(126, 13).Value = -7211.599999999999
$ws.Cells.Item(126, 14).Value = -14363

$ws.Cells.Item(132, 8).Value = 39220960
$ws.Cells.Item(132, 9).Value = 49388356
$ws.Cells.Item(132, 10).Value = 3858.8572
$ws.Cells.Item(132, 11).Value = 148165068
$ws.Cells.Item(132, 12).Value = 11576.5716
$ws.Cells.Item(132, 13).Value = -148162538
$ws.Cells.Item(132, 14).Value = -16636.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 341.2
$ws.Cells.Item(92, 10).Value = 301.5
$ws.Cells.Item(92, 12).Value = 904.5
$ws.Cells.Item(92, 14).Value = -3400.5

$ws.Cells.Item(98, 8).Value = 588.1786
$ws.Cells.Item(98, 9).Value = 415.77777
$ws.Cells.Item(98, 10).Value = 898.5
$ws.Cells.Item(98, 11).Value = 1247.33331
$ws.Cells.Item(98, 12).Value = 2695.5
$ws.Cells.Item(98, 13).Value = 250.66669
$ws.Cells.Item(98, 14).Value = -5691.5

$ws.Cells.Item(107, 8).Value = 206.48
$ws.Cells.Item(107, 9).Value = 160.88889
$ws.Cells.Item(107, 10).Value = 232.125
$ws.Cells.Item(107, 11).Value = 482.66667
$ws.Cells.Item(107, 12).Value = 696.375
$ws.Cells.Item(107, 13).Value = 1437.33333
$ws.Cells.Item(107, 14).Value = -4536.375

$ws.Cells.Item(122, 8).Value = 763.0454999999999
$ws.Cells.Item(122, 9).Value = 201.125
$ws.Cells.Item(122, 10).Value = 1084.1428
$ws.Cells.Item(122, 11).Value = 1810.125
$ws.Cells.Item(122, 12).Value = 9757.2852
$ws.Cells.Item(122, 13).Value = 639.875
$ws.Cells.Item(122, 14).Value = -14657.2852

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1300.3334
$ws.Cells.Item(122, 9).Value = 800
$ws.Cells.Item(122, 10).Value = 1550.5
$ws.Cells.Item(122, 11).Value = 2400
$ws.Cells.Item(122, 12).Value = 4651.5
$ws.Cells.Item(122, 13).Value = 50
$ws.Cells.Item(122, 14).Value = -9551.5

$ws.Cells.Item(126, 8).Value = 9525318
$ws.Cells.Item(126, 9).Value = 15874183
$ws.Cells.Item(126, 10).Value = 2021.3334
$ws.Cells.Item(126, 11).Value = 47622549
$ws.Cells.Item(126, 12).Value = 6064.0002
$ws.Cells.Item(126, 13).Value = -47620079
$ws.Cells.Item(126, 14).Value = -11004.0002

$ws.Cells.Item(132, 8).Value = 4914.1665
$ws.Cells.Item(132, 9).Value = 5789.385
$ws.Cells.Item(132, 10).Value = 2638.6
$ws.Cells.Item(132, 11).Value = 17368.155
$ws.Cells.Item(132, 12).Value = 7915.799999999999
$ws.Cells.Item(132, 13).Value = -14838.155
$ws.Cells.Item(132, 14).Value = -12975.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2115
$ws.Cells.Item(7, 9).Value = 1350
$ws.Cells.Item(7, 11).Value = 1350
$ws.Cells.Item(7, 13).Value = -1238

$ws.Cells.Item(22, 8).Value = 433.07693
$ws.Cells.Item(22, 9).Value = 349.25
$ws.Cells.Item(22, 10).Value = 470.33334
$ws.Cells.Item(22, 11).Value = 349.25
$ws.Cells.Item(22, 12).Value = 470.33334
$ws.Cells.Item(22, 13).Value = -54.25
$ws.Cells.Item(22, 14).Value = -1060.33334

$ws.Cells.Item(26, 8).Value = 10000
$ws.Cells.Item(26, 10).Value = 10000
$ws.Cells.Item(26, 12).Value = 10000
$ws.Cells.Item(26, 14).Value = -10590

$ws.Cells.Item(27, 8).Value = 433.07693
$ws.Cells.Item(27, 9).Value = 349.25
$ws.Cells.Item(27, 10).Value = 470.33334
$ws.Cells.Item(27, 11).Value = 349.25
$ws.Cells.Item(27, 12).Value = 470.33334
$ws.Cells.Item(27, 13).Value = -242.25
$ws.Cells.Item(27, 14).Value = -684.33334

$ws.Cells.Item(40, 8).Value = 2092.2354
$ws.Cells.Item(40, 9).Value = 1709
$ws.Cells.Item(40, 10).Value = 3012
$ws.Cells.Item(40, 11).Value = 1709
$ws.Cells.Item(40, 12).Value = 3012
$ws.Cells.Item(40, 13).Value = -1573
$ws.Cells.Item(40, 14).Value = -3284

$ws.Cells.Item(46, 8).Value = 5788.4
$ws.Cells.Item(46, 9).Value = 5788.4
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 11).Value = 5788.4
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).Value = -5600.4
$ws.Cells.Item(46, 14).Value = ""

$ws.Cells.Item(55, 8).Value = 164.6
$ws.Cells.Item(55, 9).Value = 60.76923
$ws.Cells.Item(55, 10).Value = 277.08334
$ws.Cells.Item(55, 11).Value = 60.76923
$ws.Cells.Item(55, 12).Value = 277.08334
$ws.Cells.Item(55, 13).Value = 112.23077
$ws.Cells.Item(55, 14).Value = -623.08334

$ws.Cells.Item(126, 8).Value = 2115
$ws.Cells.Item(126, 9).Value = 1350
$ws.Cells.Item(126, 11).Value = 4050
$ws.Cells.Item(126, 13).Value = -1580

$ws.Cells.Item(136, 8).Value = 3505.319
$ws.Cells.Item(136, 9).Value = 1375
$ws.Cells.Item(136, 10).Value = 10477.272
$ws.Cells.Item(136, 11).Value = 4125
$ws.Cells.Item(136, 12).Value = 31431.816
$ws.Cells.Item(136, 13).Value = -1575
$ws.Cells.Item(136, 14).Value = -36531.81600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 2
$ws.Cells.Item(2, 9).Value = 2
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 110
$ws.Cells.Item(2, 14).Value = ""

$ws.Cells.Item(110, 8).Value = 48000
$ws.Cells.Item(110, 10).Value = 48000
$ws.Cells.Item(110, 12).Value = 48000
$ws.Cells.Item(110, 14).Value = -56180

$ws.Cells.Item(122, 8).Value = 11112564
$ws.Cells.Item(122, 9).Value = 20000950
$ws.Cells.Item(122, 10).Value = 2079.375
$ws.Cells.Item(122, 11).Value = 60002850
$ws.Cells.Item(122, 12).Value = 6238.125
$ws.Cells.Item(122, 13).Value = -60000400
$ws.Cells.Item(122, 14).Value = -11138.125

$ws.Cells.Item(125, 8).Value = 19750
$ws.Cells.Item(125, 10).Value = 19750
$ws.Cells.Item(125, 12).Value = 19750
$ws.Cells.Item(125, 14).Value = -29590

$ws.Cells.Item(126, 8).Value = 1007.2727
$ws.Cells.Item(126, 9).Value = 645
$ws.Cells.Item(126, 10).Value = 1442
$ws.Cells.Item(126, 11).Value = 1935
$ws.Cells.Item(126, 12).Value = 4326
$ws.Cells.Item(126, 13).Value = 535
$ws.Cells.Item(126, 14).Value = -9266

$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(133, 14).Value = ""
